$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GLW")

# Row 4 - Inventory
$ws.Range("B4").Value = 2361000000.0
$ws.Range("C4").Value = 2438000000.0
$ws.Range("D4").Value = 2581000000.0
$ws.Range("E4").Value = 2235000000.0
$ws.Range("F4").Value = 2347000000.0

# Row 15 - Accounts Payable
$ws.Range("B15").Value = 1272000000.0
$ws.Range("C15").Value = 1174000000.0
$ws.Range("D15").Value = 1176000000.0
$ws.Range("E15").Value = 1109000000.0
$ws.Range("F15").Value = 1250000000.0

# Row 22 - Long Term Tax Liability (Deferred)
$ws.Range("B22").Value = -992000000.0
$ws.Range("C22").Value = -808000000.0
$ws.Range("D22").Value = -926000000.0
$ws.Range("E22").Value = -993000000.0
$ws.Range("F22").Value = -930000000.0

# Row 37 - Net Debt
$ws.Range("G37").Value = 5306000000.0

# Row 38 - Total Debt
$ws.Range("G38").Value = 7740000000.0
